$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: date in A11, count in B11 (continuing the existing table)
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 45973
$ws.Range("B11").Value = 11

# Update the active selection to match the saved view state
$ws.Range("A14").Select()
